$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tbl = $sh.Table
$tbl.FirstRow = $true
